# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to match the latest scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map: sheet name -> { row number -> new F value }
$sheetUpdates = @{
    "展览" = @{
        7  = 2683
        8  = 1165
        9  = 268
        10 = 121
        11 = 10096
        13 = 261
        15 = 626
        16 = 11772
        17 = 12139
    }
    "全部类型" = @{
        7  = 2683
        9  = 1165
        10 = 268
        11 = 121
        12 = 10096
        14 = 261
        16 = 626
        17 = 11772
        18 = 12139
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $ws.Cells.Item($rowNum, 6).Value = $rows[$rowNum]
    }
}
